$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Regenerated save_data: column G ("K") now reflects strikeouts (K) instead
# of the previous "Strike#" metric. Update the computed s_vals for each row.
$ws.Range("G2").Value = 1
$ws.Range("G3").Value = 0
$ws.Range("G5").Value = 0
$ws.Range("G6").Value = 0
$ws.Range("G7").Value = 2
$ws.Range("G8").Value = 2
